# Add a new "Test Cases" sheet at the front of the workbook, listing each
# test sheet's run mode (mirrors the existing "Keywords" / "TestA" sheets).

$wb = $excel.ActiveWorkbook

# Insert the new sheet before the current first sheet ("Keywords") so it
# becomes the new first (left-most) tab.
$beforeSheet = $wb.Worksheets.Item(1)
$testCases = $wb.Worksheets.Add($beforeSheet)
$testCases.Name = "Test Cases"

# Header row
$testCases.Range("A1").Value = "TCID"
$testCases.Range("B1").Value = "Runmode"

# One row per test-suite sheet
$testCases.Range("A2").Value = "TestA"
$testCases.Range("B2").Value = "Y"
$testCases.Range("A3").Value = "TestB"
$testCases.Range("B3").Value = "N"
$testCases.Range("A4").Value = "TestC"
$testCases.Range("B4").Value = "N"

# Match the green header-row fill used on the other sheets
$testCases.Range("A1:B1").Interior.Color = 5296274

# Widen column A so the TCID values aren't clipped
$testCases.Columns.Item(1).ColumnWidth = 14.67

# Leave the same kind of "click below the data" selection the other sheets
# have, and make this new sheet the active tab.
$null = $testCases.Range("B7").Select()
$testCases.Activate()
